# Auto-generated Excel COM-interop script to apply the Masamune_Profits data refresh
# Updates market-board derived columns (currentAveragePrice*, LevePrice*, LeveProfit*)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets to match the scheduled runner's refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 2674.95
$ws.Range("I43").Value = 2882.353
$ws.Range("J43").Value = 1499.6666
$ws.Range("K43").Value = 2882.353
$ws.Range("L43").Value = 1499.6666
$ws.Range("M43").Value = -2813.353
$ws.Range("N43").Value = -1637.6666
$ws.Range("H70").Value = 1133.1364
$ws.Range("I70").Value = 995.5625
$ws.Range("K70").Value = 2986.6875
$ws.Range("M70").Value = -2716.6875
$ws.Range("H73").Value = 1133.1364
$ws.Range("I73").Value = 995.5625
$ws.Range("K73").Value = 2986.6875
$ws.Range("M73").Value = -2050.6875
$ws.Range("H74").Value = 4646.3335
$ws.Range("I74").Value = 6500
$ws.Range("J74").Value = 3719.5
$ws.Range("K74").Value = 6500
$ws.Range("L74").Value = 3719.5
$ws.Range("M74").Value = -5564
$ws.Range("N74").Value = -5591.5
$ws.Range("I76").Value = 3000.1035
$ws.Range("J76").Value = 3100
$ws.Range("K76").Value = 3000.1035
$ws.Range("L76").Value = 3100
$ws.Range("M76").Value = -2685.1035
$ws.Range("N76").Value = -3730
$ws.Range("H77").Value = 4646.3335
$ws.Range("I77").Value = 6500
$ws.Range("J77").Value = 3719.5
$ws.Range("K77").Value = 32500
$ws.Range("L77").Value = 18597.5
$ws.Range("M77").Value = -27820
$ws.Range("N77").Value = -27957.5
$ws.Range("I79").Value = 3000.1035
$ws.Range("J79").Value = 3100
$ws.Range("K79").Value = 3000.1035
$ws.Range("L79").Value = 3100
$ws.Range("M79").Value = -1908.1035
$ws.Range("N79").Value = -5284
$ws.Range("H129").Value = 294116.28
$ws.Range("J129").Value = 1319.7222
$ws.Range("L129").Value = 3959.1666
$ws.Range("N129").Value = -13959.1666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17141.375
$ws.Range("I32").Value = 16221.152
$ws.Range("J32").Value = 28000
$ws.Range("K32").Value = 16221.152
$ws.Range("L32").Value = 28000
$ws.Range("M32").Value = -15934.152
$ws.Range("N32").Value = -28574
$ws.Range("H63").Value = 3221.6667
$ws.Range("I63").Value = 2452.5
$ws.Range("J63").Value = 4760
$ws.Range("K63").Value = 2452.5
$ws.Range("L63").Value = 4760
$ws.Range("M63").Value = -1766.5
$ws.Range("N63").Value = -6132
$ws.Range("H66").Value = 3221.6667
$ws.Range("I66").Value = 2452.5
$ws.Range("J66").Value = 4760
$ws.Range("K66").Value = 12262.5
$ws.Range("L66").Value = 23800
$ws.Range("M66").Value = -8830.5
$ws.Range("N66").Value = -30664
$ws.Range("H88").Value = 48566930
$ws.Range("I88").Value = 85717050
$ws.Range("J88").Value = 8559105
$ws.Range("K88").Value = 85717050
$ws.Range("L88").Value = 8559105
$ws.Range("M88").Value = -85716644
$ws.Range("N88").Value = -8559917
$ws.Range("H91").Value = 48566930
$ws.Range("I91").Value = 85717050
$ws.Range("J91").Value = 8559105
$ws.Range("K91").Value = 85717050
$ws.Range("L91").Value = 8559105
$ws.Range("M91").Value = -85715646
$ws.Range("N91").Value = -8561913
$ws.Range("H102").Value = 16337.143
$ws.Range("I102").Value = 2466.5
$ws.Range("J102").Value = 26740.125
$ws.Range("K102").Value = 2466.5
$ws.Range("L102").Value = 26740.125
$ws.Range("M102").Value = -844.5
$ws.Range("N102").Value = -29984.125
$ws.Range("H131").Value = 50307
$ws.Range("J131").Value = 50307
$ws.Range("L131").Value = 50307
$ws.Range("N131").Value = -60387
$ws.Range("H133").Value = 32460.73
$ws.Range("J133").Value = 32460.73
$ws.Range("L133").Value = 32460.73
$ws.Range("N133").Value = -37520.73

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2428.5715
$ws.Range("I86").Value = 2533.3333
$ws.Range("K86").Value = 2533.3333
$ws.Range("M86").Value = -1410.3333
$ws.Range("H89").Value = 2428.5715
$ws.Range("I89").Value = 2533.3333
$ws.Range("K89").Value = 12666.6665
$ws.Range("M89").Value = -7050.666499999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2758.8948
$ws.Range("I62").Value = 2608.5
$ws.Range("J62").Value = 3180
$ws.Range("K62").Value = 2608.5
$ws.Range("L62").Value = 3180
$ws.Range("M62").Value = -1984.5
$ws.Range("N62").Value = -4428
$ws.Range("H65").Value = 2758.8948
$ws.Range("I65").Value = 2608.5
$ws.Range("J65").Value = 3180
$ws.Range("K65").Value = 13042.5
$ws.Range("L65").Value = 15900
$ws.Range("M65").Value = -9922.5
$ws.Range("N65").Value = -22140
$ws.Range("H99").Value = 1886.7222
$ws.Range("I99").Value = 1946.4546
$ws.Range("J99").Value = 1792.8572
$ws.Range("K99").Value = 1946.4546
$ws.Range("L99").Value = 1792.8572
$ws.Range("M99").Value = -448.4546
$ws.Range("N99").Value = -4788.8572
$ws.Range("H126").Value = 1886.7222
$ws.Range("I126").Value = 1946.4546
$ws.Range("J126").Value = 1792.8572
$ws.Range("K126").Value = 5839.3638
$ws.Range("L126").Value = 5378.571599999999
$ws.Range("M126").Value = -3369.3638
$ws.Range("N126").Value = -10318.5716

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H101").Value = 5999.3335
$ws.Range("J101").Value = 5999.3335
$ws.Range("L101").Value = 17998.0005
$ws.Range("N101").Value = -22866.0005
$ws.Range("H102").Value = 14124.75
$ws.Range("J102").Value = 20000
$ws.Range("L102").Value = 60000
$ws.Range("N102").Value = -64868
$ws.Range("H106").Value = 1303014.5
$ws.Range("J106").Value = 6029
$ws.Range("L106").Value = 18087
$ws.Range("N106").Value = -19979

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5313.378
$ws.Range("I70").Value = 5188.9443
$ws.Range("J70").Value = 5811.1113
$ws.Range("K70").Value = 5188.9443
$ws.Range("L70").Value = 5811.1113
$ws.Range("M70").Value = -4918.9443
$ws.Range("N70").Value = -6351.1113
$ws.Range("H73").Value = 5313.378
$ws.Range("I73").Value = 5188.9443
$ws.Range("J73").Value = 5811.1113
$ws.Range("K73").Value = 5188.9443
$ws.Range("L73").Value = 5811.1113
$ws.Range("M73").Value = -4252.9443
$ws.Range("N73").Value = -7683.1113
$ws.Range("H80").Value = 361428.94
$ws.Range("I80").Value = 504630.5
$ws.Range("J80").Value = 3425
$ws.Range("K80").Value = 504630.5
$ws.Range("L80").Value = 3425
$ws.Range("M80").Value = -503632.5
$ws.Range("N80").Value = -5421
$ws.Range("H83").Value = 361428.94
$ws.Range("I83").Value = 504630.5
$ws.Range("J83").Value = 3425
$ws.Range("K83").Value = 2523152.5
$ws.Range("L83").Value = 17125
$ws.Range("M83").Value = -2518160.5
$ws.Range("N83").Value = -27109
$ws.Range("H97").Value = 11221.5
$ws.Range("I97").Value = 3805
$ws.Range("J97").Value = 18638
$ws.Range("K97").Value = 3805
$ws.Range("L97").Value = 18638
$ws.Range("M97").Value = -3309
$ws.Range("N97").Value = -19630
$ws.Range("H113").Value = 1088.25
$ws.Range("I113").Value = 1025.4166
$ws.Range("J113").Value = 1182.5
$ws.Range("K113").Value = 1025.4166
$ws.Range("L113").Value = 1182.5
$ws.Range("M113").Value = 1144.5834
$ws.Range("N113").Value = -5522.5
$ws.Range("H122").Value = 1756.25
$ws.Range("I122").Value = 1966.6666
$ws.Range("J122").Value = 1630
$ws.Range("K122").Value = 5899.9998
$ws.Range("L122").Value = 4890
$ws.Range("M122").Value = -3449.9998
$ws.Range("N122").Value = -9790

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1841.381
$ws.Range("I46").Value = 1342.8572
$ws.Range("J46").Value = 2090.6428
$ws.Range("K46").Value = 1342.8572
$ws.Range("L46").Value = 2090.6428
$ws.Range("M46").Value = -1154.8572
$ws.Range("N46").Value = -2466.6428
$ws.Range("H68").Value = 2214.6428
$ws.Range("I68").Value = 2257.7144
$ws.Range("J68").Value = 2171.5715
$ws.Range("K68").Value = 2257.7144
$ws.Range("L68").Value = 2171.5715
$ws.Range("M68").Value = -1508.7144
$ws.Range("N68").Value = -3669.5715
$ws.Range("H71").Value = 2214.6428
$ws.Range("I71").Value = 2257.7144
$ws.Range("J71").Value = 2171.5715
$ws.Range("K71").Value = 11288.572
$ws.Range("L71").Value = 10857.8575
$ws.Range("M71").Value = -7544.572
$ws.Range("N71").Value = -18345.8575
$ws.Range("H122").Value = 44105.457
$ws.Range("I122").Value = 64408.188
$ws.Range("K122").Value = 193224.564
$ws.Range("M122").Value = -190774.564

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 74842.875
$ws.Range("J46").Value = 82479
$ws.Range("L46").Value = 82479
$ws.Range("N46").Value = -82941
$ws.Range("H62").Value = 3000
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 3000
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 3000
$ws.Range("N62").Value = -4248
$ws.Range("H65").Value = 3000
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 3000
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 15000
$ws.Range("N65").Value = -21240
$ws.Range("H122").Value = 75511130
$ws.Range("I122").Value = 88096150
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 264288450
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -264286000
$ws.Range("N122").Value = -7900
$ws.Range("H123").Value = 38336
$ws.Range("J123").Value = 38336
$ws.Range("L123").Value = 38336
$ws.Range("N123").Value = -48136
$ws.Range("H134").Value = 74842.875
$ws.Range("J134").Value = 82479
$ws.Range("L134").Value = 247437
$ws.Range("N134").Value = -252507
$ws.Range("M62").ClearContents()
$ws.Range("M65").ClearContents()
